$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.917.24"
$ws.Range("E2").Value = "  +1.13%  "

$ws.Range("D3").Value = "'2.489.63"

$ws.Range("D5").Value = "'588.55"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").Value = "'174.43"
$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +3.61%  "

$ws.Range("E10").Value = "  -1.48%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").Value = "'67.868.84"
$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("D17").Value = "'2.521.65"
$ws.Range("E17").Value = "  +3.01%  "

$ws.Range("D18").Value = "'10.86"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D20").Value = "'347.29"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("E21").Value = "  +2.38%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D24").Value = "'4.17"
$ws.Range("E24").Value = "  -1.59%  "

$ws.Range("E25").Value = "  -5.48%  "

$ws.Range("D26").Value = "'8.84"
$ws.Range("E26").Value = "  -3.05%  "

$ws.Range("E27").Value = "  +0.32%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D29").Value = "'0.0₃0890"
$ws.Range("E29").Value = "  -2.01%  "

$ws.Range("D30").Value = "'497.60"
$ws.Range("E30").Value = "  -1.30%  "

$ws.Range("D31").Value = "'7.75"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").Value = "'164.12"
$ws.Range("E35").Value = "  +0.94%  "

$ws.Range("E36").Value = "  +1.99%  "

$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("D38").Value = "'18.22"
$ws.Range("E38").Value = "  +0.33%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  -1.71%  "

$ws.Range("E41").Value = "  +2.61%  "

$ws.Range("D43").Value = "'4.78"
$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("D44").Value = "'2.38"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").Value = "'148.48"
$ws.Range("E45").Value = "  +3.77%  "

$ws.Range("E46").Value = "  +1.47%  "

$ws.Range("D47").Value = "'0.513"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("E48").Value = "  -4.62%  "

$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("E51").Value = "  -1.38%  "
